# Testing out new ig pub approach with templates.
# The "valueCode" slice row (old row 9) and the "valueString" slice row
# (old row 16) are removed from the Elements table; everything below each
# shifts up. A handful of cells are also retouched once the rows have
# settled into their final positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete slice rows (delete the lower one first so the
# second deletion index still refers to the intended row).
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(9).Delete()

# --- Cell-level touch-ups on the now-shifted table ------------------------

# Row 8 (Extension.extension.url / code slice): Type(s) becomes "string".
$ws.Range("J8").Value = "string`n"

# Row 9 (Extension.extension.value[x] / code slice): no longer has its own
# Slice Name (that lived on the row that got removed above).
$ws.Range("B9").Value = ""

# Row 13 (Extension.extension.url / text slice): Type(s) becomes "string".
$ws.Range("J13").Value = "string`n"

# Row 14 (Extension.extension.value[x] / text slice): clear the leftover
# slicing metadata that belonged to the deleted row.
$ws.Range("AA14").Value = ""
$ws.Range("AB14").Value = ""
$ws.Range("AD14").Value = ""

# Row 15 (Extension.url): now typed as "uri" with a fixed value pointing at
# this extension's own StructureDefinition.
$ws.Range("J15").Value = "uri`n"
$ws.Range("Q15").Value = "http://www.fhir.org/guides/test3/StructureDefinition/extension-complex"

# Row 16 (Extension.value[x]): Max cardinality drops to 0.
$ws.Range("F16").Value = "0"

# --- Re-sync the structural metadata with the new 16-row extent ----------

# Defined name backing the (hidden) filter database.
$wb.Names.Item(1).RefersTo = "=Elements!`$A`$1:`$AJ`$16"

# AutoFilter: drop and reapply over the shrunk range with the same two
# column filters it had before.
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ16").AutoFilter(7, "<> ")
$ws.Range("A1:AJ16").AutoFilter(27, @(""))

# Conditional formatting previously covered A2:AI17; now covers A2:AI15.
$fcs = $ws.Cells.Item(2, 1).FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2:AI15"))
}
